# Commit: "updated sierpinski benchmark scene"
# Inserts a new "Mandelbulb Optimisations" worksheet between "Mandelbulb" and
# "Sierpinski", populated with a benchmark table (same shape as the existing
# "Mandelbulb" sheet), and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet right after "Mandelbulb" (i.e. right before
#    "Sierpinski"), and name it.
# ---------------------------------------------------------------------------
$mandelbulb = $wb.Worksheets.Item("Mandelbulb")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mandelbulb)
$ws.Name = "Mandelbulb Optimisations"

# ---------------------------------------------------------------------------
# 2. Column widths (matches the other benchmark sheets' layout).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.42578125
$ws.Columns.Item(2).ColumnWidth = 54
$ws.Columns.Item(3).ColumnWidth = 95
$ws.Columns.Item(4).ColumnWidth = 7.140625
$ws.Columns.Item(5).ColumnWidth = 7.140625
$ws.Columns.Item(6).ColumnWidth = 26.5703125
$ws.Columns.Item(7).ColumnWidth = 16.42578125
$ws.Columns.Item(8).ColumnWidth = 16.28515625
$ws.Columns.Item(9).ColumnWidth = 21.85546875
$ws.Columns.Item(10).ColumnWidth = 22.85546875
$ws.Columns.Item(11).ColumnWidth = 22.28515625
$ws.Columns.Item(12).ColumnWidth = 20.7109375
$ws.Columns.Item(13).ColumnWidth = 24.140625
$ws.Columns.Item(14).ColumnWidth = 17.5703125
$ws.Columns.Item(15).ColumnWidth = 13.140625
$ws.Columns.Item(16).ColumnWidth = 24.7109375
$ws.Columns.Item(17).ColumnWidth = 24.5703125
$ws.Columns.Item(18).ColumnWidth = 22.28515625
$ws.Columns.Item(19).ColumnWidth = 22.28515625
$ws.Columns.Item(20).ColumnWidth = 19.85546875
$ws.Columns.Item(21).ColumnWidth = 18
$ws.Columns.Item(22).ColumnWidth = 21.28515625
$ws.Columns.Item(23).ColumnWidth = 23
$ws.Columns.Item(24).ColumnWidth = 12.28515625
$ws.Columns.Item(25).ColumnWidth = 13.7109375
$ws.Columns.Item(26).ColumnWidth = 13.5703125
$ws.Columns.Item(27).ColumnWidth = 11.7109375
$ws.Columns.Item(28).ColumnWidth = 11.5703125

# ---------------------------------------------------------------------------
# 3. Header row (row 1) - green fill style (same "Style2" used elsewhere).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Scene Description"
$ws.Range("B1").Value = "Scene name"
$ws.Range("C1").Value = "Build Options"
$ws.Range("D1").Value = "Width"
$ws.Range("E1").Value = "Height"
$ws.Range("F1").Value = "Device Name"
$ws.Range("G1").Value = "Device Version"
$ws.Range("H1").Value = "Work Group Size"
$ws.Range("I1").Value = "Clock Frequency (MHz)"
$ws.Range("J1").Value = "Parallel Compute Units"
$ws.Range("K1").Value = "Global Memory (Bytes)"
$ws.Range("L1").Value = "Local Memory (Bytes)"
$ws.Range("M1").Value = "Constant Memory (Bytes)"
$ws.Range("N1").Value = "Total Duration (s)"
$ws.Range("O1").Value = "Total Frames"
$ws.Range("P1").Value = "Maximum Frame Time (s)"
$ws.Range("Q1").Value = "Minimum Frame Time (s)"
$ws.Range("A1:Q1").Style = $mandelbulb.Range("A1:Q1").Style

$ws.Range("S1").Value = "Total Number of Pixels"
$ws.Range("T1").Value = "Global Memory (GB)"
$ws.Range("U1").Value = "Local Memory (KB)"
$ws.Range("V1").Value = "Constant Memory (KB)"
$ws.Range("W1").Value = "Mean Frame Time (s)"
$ws.Range("X1").Value = "Mean FPS"
$ws.Range("Y1").Value = "Maximum FPS"
$ws.Range("Z1").Value = "Minimum FPS"
$ws.Range("AA1").Value = "Max - Mean"
$ws.Range("AB1").Value = "Mean - Min"
$ws.Range("S1:AB1").Style = $mandelbulb.Range("S1:AB1").Style

# ---------------------------------------------------------------------------
# 4. Data rows. Values are typed in the same order the original author
#    would have entered them, so brand-new shared strings land at the same
#    shared-string-table offsets as the authored workbook (Bounding Volume,
#    Linear Epsilon, None, All, "-cl-fast-relaxed-math ").
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Bounding Volume"
$ws.Range("A4").Value = "Linear Epsilon"
$ws.Range("A2").Value = "None"
$ws.Range("A6").Value = "All"
$ws.Range("A5").Value = "-cl-fast-relaxed-math "

$ws.Range("B2").Value = "kernels/benchmarks/mandelbulb_base.cl"
$ws.Range("B3").Value = "kernels/benchmarks/mandelbulb_bounding_volume.cl"
$ws.Range("B4").Value = "kernels/benchmarks/mandelbulb_linear_epsilon.cl"
$ws.Range("B5").Value = "kernels/benchmarks/mandelbulb_base.cl"
$ws.Range("B6").Value = "kernels/benchmarks/mandelbulb_ALL.cl"

$ws.Range("C2").Value = '-I "kernels" -I "kernels\benchmarks" -I "kernels\include" -I "kernels/include" '
$ws.Range("C3").Value = '-I "kernels" -I "kernels\benchmarks" -I "kernels\include" -I "kernels/include" '
$ws.Range("C4").Value = '-I "kernels" -I "kernels\benchmarks" -I "kernels\include" -I "kernels/include" '
$ws.Range("C5").Value = '-I "kernels" -I "kernels\benchmarks" -I "kernels\include" -I "kernels/include" -cl-fast-relaxed-math'
$ws.Range("C6").Value = '-I "kernels" -I "kernels\benchmarks" -I "kernels\include" -I "kernels/include" -cl-fast-relaxed-math'

for ($r = 2; $r -le 6; $r++) {
    $ws.Range("D$r").Value = 1920
    $ws.Range("E$r").Value = 1080
    $ws.Range("F$r").Value = "NVIDIA GeForce RTX 3060 Ti"
    $ws.Range("G$r").Value = "OpenCL 3.0 CUDA"
    $ws.Range("H$r").Value = 256
    $ws.Range("I$r").Value = 1665
    $ws.Range("J$r").Value = 38
    $ws.Range("K$r").Value = 8589410304
    $ws.Range("L$r").Value = 49152
    $ws.Range("M$r").Value = 65536
}

$ws.Range("N2").Value = 39.0451
$ws.Range("O2").Value = 502
$ws.Range("P2").Value = 0.101914
$ws.Range("Q2").Value = 0.0512617

$ws.Range("N3").Value = 39.0797
$ws.Range("O3").Value = 468
$ws.Range("P3").Value = 0.111606
$ws.Range("Q3").Value = 0.0537429

$ws.Range("N4").Value = 39.1144
$ws.Range("O4").Value = 623
$ws.Range("P4").Value = 0.094061
$ws.Range("Q4").Value = 0.033296

$ws.Range("N5").Value = 39.0618
$ws.Range("O5").Value = 993
$ws.Range("P5").Value = 0.0538985
$ws.Range("Q5").Value = 0.0275404

$ws.Range("N6").Value = 39.0419
$ws.Range("O6").Value = 1225
$ws.Range("P6").Value = 0.0447649
$ws.Range("Q6").Value = 0.0188609

$ws.Range("A2:Q6").Style = $mandelbulb.Range("A2:Q6").Style

# ---------------------------------------------------------------------------
# 5. Derived / helper columns (S:AB). Rows 2-4 are "normal"; rows 5 & 6
#    cross-reference each other for S,T,U,V,Y,Z (copy/paste artefact carried
#    over from the source workbook) while X/AA/AB stay self-referential.
# ---------------------------------------------------------------------------
$ws.Range("S2").Formula = "=D2*E2"
$ws.Range("T2").Formula = "=K2/1000000000"
$ws.Range("U2").Formula = "=L2/1000"
$ws.Range("V2").Formula = "=M2/1000"
$ws.Range("W2").Formula = "=N2/O2"
$ws.Range("X2").Formula = "=1/W2"
$ws.Range("Y2").Formula = "=1/Q2"
$ws.Range("Z2").Formula = "=1/P2"
$ws.Range("AA2").Formula = "=Y2-X2"
$ws.Range("AB2").Formula = "=X2-Z2"

$ws.Range("S3").Formula = "=D3*E3"
$ws.Range("T3").Formula = "=K3/1000000000"
$ws.Range("U3").Formula = "=L3/1000"
$ws.Range("V3").Formula = "=M3/1000"
$ws.Range("W3").Formula = "=N3/O3"
$ws.Range("X3").Formula = "=1/W3"
$ws.Range("Y3").Formula = "=1/Q3"
$ws.Range("Z3").Formula = "=1/P3"
$ws.Range("AA3").Formula = "=Y3-X3"
$ws.Range("AB3").Formula = "=X3-Z3"

$ws.Range("S4").Formula = "=D4*E4"
$ws.Range("T4").Formula = "=K4/1000000000"
$ws.Range("U4").Formula = "=L4/1000"
$ws.Range("V4").Formula = "=M4/1000"
$ws.Range("W4").Formula = "=N4/O4"
$ws.Range("X4").Formula = "=1/W4"
$ws.Range("Y4").Formula = "=1/Q4"
$ws.Range("Z4").Formula = "=1/P4"
$ws.Range("AA4").Formula = "=Y4-X4"
$ws.Range("AB4").Formula = "=X4-Z4"

$ws.Range("S5").Formula = "=D6*E6"
$ws.Range("T5").Formula = "=K6/1000000000"
$ws.Range("U5").Formula = "=L6/1000"
$ws.Range("V5").Formula = "=M6/1000"
$ws.Range("W5").Formula = "=N6/O6"
$ws.Range("X5").Formula = "=1/W5"
$ws.Range("Y5").Formula = "=1/Q6"
$ws.Range("Z5").Formula = "=1/P6"
$ws.Range("AA5").Formula = "=Y5-X5"
$ws.Range("AB5").Formula = "=X5-Z5"

$ws.Range("S6").Formula = "=D5*E5"
$ws.Range("T6").Formula = "=K5/1000000000"
$ws.Range("U6").Formula = "=L5/1000"
$ws.Range("V6").Formula = "=M5/1000"
$ws.Range("W6").Formula = "=N5/O5"
$ws.Range("X6").Formula = "=1/W6"
$ws.Range("Y6").Formula = "=1/Q5"
$ws.Range("Z6").Formula = "=1/P5"
$ws.Range("AA6").Formula = "=Y6-X6"
$ws.Range("AB6").Formula = "=X6-Z6"

$ws.Range("T2:V6").Style = $mandelbulb.Range("T2:V6").Style

# Rows 7 & 8 only carry the T/U/V number-format style, no data/formula.
$ws.Range("T7:V8").Style = $mandelbulb.Range("T2:V2").Style

# ---------------------------------------------------------------------------
# 6. Sheet view state - matches the authored file (this sheet becomes the
#    active / selected tab, scrolled so column D is the left-most visible
#    column, with a big "select all" style range selected).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Range("A1:XFD1048576").Select()
[void]$ws.Range("S20").Select()

[void]$wb.Worksheets.Item("Graphs").Activate()
[void]$ws.Activate()
